# Auto-update draw results: append the new day's Pick 4 row.
# Mirrors the daily bot commit that appends one row (A:E) to the Results
# sheet and grows the used range from A1:E16 to A1:E17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = $ws.UsedRange.Rows.Count + 1   # 16 -> 17, next blank row

$rng = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 5))

# Every existing row stores its values as plain text (dates, the numeric
# "phase" code, and the timestamp are all text, not real numbers/dates), so
# force Text format on the new row before assigning values - otherwise
# Excel would auto-coerce "2025-10-03" into a date serial and "251003"
# into a number.
$rng.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-10-03"
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Value = "251003"
$ws.Cells.Item($row, 4).Value = "2-8-2-1"
$ws.Cells.Item($row, 5).Value = "2025-10-03T21:36:35.491+04:00"
